$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (财联社) and column C (同花顺), rows 2-21.
# Column B (东方财富) is cleared for all these rows.
$colA = @(
    "利欧股份",
    "华胜天成",
    "大位科技",
    "博纳影业",
    "掌阅科技",
    "光线传媒",
    "天奇股份",
    "巨力索具",
    "协鑫集成",
    "深科技",
    "横店影视",
    "网宿科技",
    "海兰信",
    "汉缆股份",
    "双良节能",
    "航发动力",
    "风语筑",
    "五洲新春",
    "万向钱潮",
    "优刻得-W"
)

$colC = @(
    "博纳影业",
    "华胜天成",
    "大位科技",
    "巨力索具",
    "光线传媒",
    "掌阅科技",
    "利欧股份",
    "天奇股份",
    "嘉美包装",
    "协鑫集成",
    "汉缆股份",
    "深科技",
    "横店影视",
    "风语筑",
    "航发动力",
    "特发信息",
    "天汽模",
    "万向钱潮",
    "浙江世宝",
    "双良节能"
)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = ""
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}
